# Insert a new daily price record for "Cebollín" (row 146) into the
# "Vega Modelo de Temuco" sheet, pushing every existing row from 146
# downward down by one (old row 146 -> 147, ..., old row 208 -> 209).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 146; everything below shifts down.
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new record's data.
$ws.Range("A146").Value = 10
$ws.Range("B146").Value = "Vega Modelo de Temuco"
$ws.Range("C146").Value = "La Araucanía"
$ws.Range("D146").Value = 44466
$ws.Range("E146").Value = 9
$ws.Range("F146").Value = 100112037
$ws.Range("G146").Value = "Cebollín"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 50
$ws.Range("K146").Value = 5000
$ws.Range("L146").Value = 5000
$ws.Range("M146").Value = 5000
$ws.Range("N146").Value = "$/docena de paquetes"
$ws.Range("O146").Value = "Región de O'Higgins"
$ws.Range("P146").Value = 417
$ws.Range("Q146").Value = 12
$ws.Range("R146").Value = "Hortaliza"
